$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 210,7
$data[0,0] = 0
$data[0,1] = 'MADRI-MADR-H-03-COS'
$data[0,2] = '100:0'
$data[0,3] = 'unlocked'
$data[0,4] = 'MADRI-MADR-H-01-DAAS'
$data[0,5] = 'xe-0/0/14'
$data[0,6] = 'PUERTOLIBRE'
$data[1,0] = 1
$data[1,1] = 'MADRI-MADR-H-03-COS'
$data[1,2] = '100:1'
$data[1,3] = 'unlocked'
$data[1,4] = 'MADRI-MADR-H-01-DAAS'
$data[1,5] = 'xe-0/0/15'
$data[1,6] = 'PUERTOLIBRE'
$data[2,0] = 2
$data[2,1] = 'MADRI-MADR-H-03-COS'
$data[2,2] = '101:0'
$data[2,3] = 'unlocked'
$data[2,4] = 'MADRI-MADR-H-01-DAAS'
$data[2,5] = 'xe-0/0/16'
$data[2,6] = 'PUERTOLIBRE'
$data[3,0] = 3
$data[3,1] = 'MADRI-MADR-H-03-COS'
$data[3,2] = '101:1'
$data[3,3] = 'unlocked'
$data[3,4] = 'MADRI-MADR-H-01-DAAS'
$data[3,5] = 'xe-0/0/17'
$data[3,6] = 'PUERTOLIBRE'
$data[4,0] = 4
$data[4,1] = 'MADRI-MADR-H-03-COS'
$data[4,2] = '102:0'
$data[4,3] = 'unlocked'
$data[4,4] = 'MADRI-MADR-H-01-DAAS'
$data[4,5] = 'xe-0/0/18'
$data[4,6] = 'PUERTOLIBRE'
$data[5,0] = 5
$data[5,1] = 'MADRI-MADR-H-03-COS'
$data[5,2] = '102:1'
$data[5,3] = 'unlocked'
$data[5,4] = 'MADRI-MADR-H-01-DAAS'
$data[5,5] = 'xe-0/0/19'
$data[5,6] = 'PUERTOLIBRE'
$data[6,0] = 6
$data[6,1] = 'MADRI-MADR-H-03-COS'
$data[6,2] = '103:0'
$data[6,3] = 'unlocked'
$data[6,4] = 'MADRI-MADR-H-01-DAAS'
$data[6,5] = 'xe-0/0/20'
$data[6,6] = 'PUERTOLIBRE'
$data[7,0] = 7
$data[7,1] = 'MADRI-MADR-H-03-COS'
$data[7,2] = '103:1'
$data[7,3] = 'unlocked'
$data[7,4] = 'MADRI-MADR-H-01-DAAS'
$data[7,5] = 'xe-0/0/21'
$data[7,6] = 'PUERTOLIBRE'
$data[8,0] = 8
$data[8,1] = 'MADRI-MADR-H-03-COS'
$data[8,2] = '104:0'
$data[8,3] = 'unlocked'
$data[8,4] = 'MADRI-MADR-H-01-DAAS'
$data[8,5] = 'xe-0/0/22'
$data[8,6] = 'PUERTOLIBRE'
$data[9,0] = 9
$data[9,1] = 'MADRI-MADR-H-03-COS'
$data[9,2] = '104:1'
$data[9,3] = 'unlocked'
$data[9,4] = 'MADRI-MADR-H-01-DAAS'
$data[9,5] = 'xe-0/0/23'
$data[9,6] = 'PUERTOLIBRE'
$data[10,0] = 10
$data[10,1] = 'MADRI-MADR-H-03-COS'
$data[10,2] = '105:0'
$data[10,3] = 'unlocked'
$data[10,4] = 'MADRI-MADR-H-01-DAAS'
$data[10,5] = 'xe-0/0/24'
$data[10,6] = 'PUERTOLIBRE'
$data[11,0] = 11
$data[11,1] = 'MADRI-MADR-H-03-COS'
$data[11,2] = '105:1'
$data[11,3] = 'unlocked'
$data[11,4] = 'MADRI-MADR-H-01-DAAS'
$data[11,5] = 'xe-0/0/25'
$data[11,6] = 'PUERTOLIBRE'
$data[12,0] = 12
$data[12,1] = 'MADRI-MADR-H-03-COS'
$data[12,2] = '106:0'
$data[12,3] = 'unlocked'
$data[12,4] = 'MADRI-MADR-H-01-DAAS'
$data[12,5] = 'xe-0/0/26'
$data[12,6] = 'PUERTOLIBRE'
$data[13,0] = 13
$data[13,1] = 'MADRI-MADR-H-03-COS'
$data[13,2] = '106:1'
$data[13,3] = 'unlocked'
$data[13,4] = 'MADRI-MADR-H-01-DAAS'
$data[13,5] = 'xe-0/0/27'
$data[13,6] = 'PUERTOLIBRE'
$data[14,0] = 14
$data[14,1] = 'MADRI-MADR-H-03-COS'
$data[14,2] = '107:0'
$data[14,3] = 'unlocked'
$data[14,4] = 'MADRI-MADR-H-01-DAAS'
$data[14,5] = 'xe-0/0/28'
$data[14,6] = 'PUERTOLIBRE'
$data[15,0] = 15
$data[15,1] = 'MADRI-MADR-H-03-COS'
$data[15,2] = '107:1'
$data[15,3] = 'unlocked'
$data[15,4] = 'MADRI-MADR-H-01-DAAS'
$data[15,5] = 'xe-0/0/29'
$data[15,6] = 'PUERTOLIBRE'
$data[16,0] = 16
$data[16,1] = 'MADRI-MADR-H-03-COS'
$data[16,2] = '108:0'
$data[16,3] = 'unlocked'
$data[16,4] = 'MADRI-MADR-H-01-DAAS'
$data[16,5] = 'xe-0/0/30'
$data[16,6] = 'PUERTOLIBRE'
$data[17,0] = 17
$data[17,1] = 'MADRI-MADR-H-03-COS'
$data[17,2] = '108:1'
$data[17,3] = 'unlocked'
$data[17,4] = 'MADRI-MADR-H-01-DAAS'
$data[17,5] = 'xe-0/0/31'
$data[17,6] = 'PUERTOLIBRE'
$data[18,0] = 18
$data[18,1] = 'MADRI-MADR-H-03-COS'
$data[18,2] = '109:0'
$data[18,3] = 'unlocked'
$data[18,4] = 'MADRI-MADR-H-01-DAAS'
$data[18,5] = 'xe-0/0/32'
$data[18,6] = 'PUERTOLIBRE'
$data[19,0] = 19
$data[19,1] = 'MADRI-MADR-H-03-COS'
$data[19,2] = '109:1'
$data[19,3] = 'unlocked'
$data[19,4] = 'MADRI-MADR-H-01-DAAS'
$data[19,5] = 'xe-0/0/33'
$data[19,6] = 'PUERTOLIBRE'
$data[20,0] = 20
$data[20,1] = 'MADRI-MADR-H-03-COS'
$data[20,2] = '10:0'
$data[20,3] = 'unlocked'
$data[20,4] = 'MADRI-MADR-H-01-DAAS'
$data[20,5] = 'xe-0/0/34'
$data[20,6] = 'PUERTOLIBRE'
$data[21,0] = 21
$data[21,1] = 'MADRI-MADR-H-03-COS'
$data[21,2] = '10:1'
$data[21,3] = 'unlocked'
$data[21,4] = 'MADRI-MADR-H-01-DAAS'
$data[21,5] = 'xe-0/0/36'
$data[21,6] = 'PUERTOLIBRE'
$data[22,0] = 22
$data[22,1] = 'MADRI-MADR-H-03-COS'
$data[22,2] = '110:0'
$data[22,3] = 'unlocked'
$data[22,4] = 'MADRI-MADR-H-01-DAAS'
$data[22,5] = 'xe-0/0/37'
$data[22,6] = 'PUERTOLIBRE'
$data[23,0] = 23
$data[23,1] = 'MADRI-MADR-H-03-COS'
$data[23,2] = '110:1'
$data[23,3] = 'unlocked'
$data[23,4] = 'MADRI-MADR-H-01-DAAS'
$data[23,5] = 'xe-0/0/38'
$data[23,6] = 'PUERTOLIBRE'
$data[24,0] = 24
$data[24,1] = 'MADRI-MADR-H-03-COS'
$data[24,2] = '111:0'
$data[24,3] = 'unlocked'
$data[24,4] = 'MADRI-MADR-H-01-DAAS'
$data[24,5] = 'xe-0/0/39'
$data[24,6] = 'PUERTOLIBRE'
$data[25,0] = 25
$data[25,1] = 'MADRI-MADR-H-03-COS'
$data[25,2] = '111:1'
$data[25,3] = 'unlocked'
$data[25,4] = 'MADRI-MADR-H-01-DAAS'
$data[25,5] = 'xe-0/0/40'
$data[25,6] = 'PUERTOLIBRE'
$data[26,0] = 26
$data[26,1] = 'MADRI-MADR-H-03-COS'
$data[26,2] = '112:0'
$data[26,3] = 'unlocked'
$data[26,4] = 'MADRI-MADR-H-01-DAAS'
$data[26,5] = 'xe-0/0/41'
$data[26,6] = 'PUERTOLIBRE'
$data[27,0] = 27
$data[27,1] = 'MADRI-MADR-H-03-COS'
$data[27,2] = '112:1'
$data[27,3] = 'unlocked'
$data[27,4] = 'MADRI-MADR-H-01-DAAS'
$data[27,5] = 'xe-0/0/42'
$data[27,6] = 'PUERTOLIBRE'
$data[28,0] = 28
$data[28,1] = 'MADRI-MADR-H-03-COS'
$data[28,2] = '14:0'
$data[28,3] = 'unlocked'
$data[28,4] = 'MADRI-MADR-H-01-DAAS'
$data[28,5] = 'xe-0/0/43'
$data[28,6] = 'PUERTOLIBRE'
$data[29,0] = 29
$data[29,1] = 'MADRI-MADR-H-03-COS'
$data[29,2] = '14:1'
$data[29,3] = 'unlocked'
$data[29,4] = 'MADRI-MADR-H-01-DAAS'
$data[29,5] = 'xe-0/0/44'
$data[29,6] = 'PUERTOLIBRE'
$data[30,0] = 30
$data[30,1] = 'MADRI-MADR-H-03-COS'
$data[30,2] = '15:0'
$data[30,3] = 'unlocked'
$data[30,4] = 'MADRI-MADR-H-01-DAAS'
$data[30,5] = 'xe-0/0/45'
$data[30,6] = 'PUERTOLIBRE'
$data[31,0] = 31
$data[31,1] = 'MADRI-MADR-H-03-COS'
$data[31,2] = '15:1'
$data[31,3] = 'unlocked'
$data[31,4] = 'MADRI-MADR-H-01-DAAS'
$data[31,5] = 'xe-0/0/46'
$data[31,6] = 'PUERTOLIBRE'
$data[32,0] = 32
$data[32,1] = 'MADRI-MADR-H-03-COS'
$data[32,2] = '16:0'
$data[32,3] = 'unlocked'
$data[32,4] = 'MADRI-MADR-H-01-DAAS'
$data[32,5] = 'xe-0/0/47'
$data[32,6] = 'PUERTOLIBRE'
$data[33,0] = 33
$data[33,1] = 'MADRI-MADR-H-03-COS'
$data[33,2] = '16:1'
$data[33,3] = 'unlocked'
$data[33,4] = 'MADRI-MADR-H-01-DAAS'
$data[33,5] = 'xe-0/0/48'
$data[33,6] = 'PUERTOLIBRE'
$data[34,0] = 34
$data[34,1] = 'MADRI-MADR-H-03-COS'
$data[34,2] = '17:0'
$data[34,3] = 'unlocked'
$data[34,4] = ""
$data[34,5] = ""
$data[34,6] = ""
$data[35,0] = 35
$data[35,1] = 'MADRI-MADR-H-03-COS'
$data[35,2] = '17:1'
$data[35,3] = 'unlocked'
$data[35,4] = ""
$data[35,5] = ""
$data[35,6] = ""
$data[36,0] = 36
$data[36,1] = 'MADRI-MADR-H-03-COS'
$data[36,2] = '18:0'
$data[36,3] = 'unlocked'
$data[36,4] = ""
$data[36,5] = ""
$data[36,6] = ""
$data[37,0] = 37
$data[37,1] = 'MADRI-MADR-H-03-COS'
$data[37,2] = '18:1'
$data[37,3] = 'unlocked'
$data[37,4] = ""
$data[37,5] = ""
$data[37,6] = ""
$data[38,0] = 38
$data[38,1] = 'MADRI-MADR-H-03-COS'
$data[38,2] = '19:0'
$data[38,3] = 'unlocked'
$data[38,4] = ""
$data[38,5] = ""
$data[38,6] = ""
$data[39,0] = 39
$data[39,1] = 'MADRI-MADR-H-03-COS'
$data[39,2] = '19:1'
$data[39,3] = 'unlocked'
$data[39,4] = ""
$data[39,5] = ""
$data[39,6] = ""
$data[40,0] = 40
$data[40,1] = 'MADRI-MADR-H-03-COS'
$data[40,2] = '20:0'
$data[40,3] = 'unlocked'
$data[40,4] = ""
$data[40,5] = ""
$data[40,6] = ""
$data[41,0] = 41
$data[41,1] = 'MADRI-MADR-H-03-COS'
$data[41,2] = '20:1'
$data[41,3] = 'unlocked'
$data[41,4] = ""
$data[41,5] = ""
$data[41,6] = ""
$data[42,0] = 42
$data[42,1] = 'MADRI-MADR-H-03-COS'
$data[42,2] = '21:0'
$data[42,3] = 'unlocked'
$data[42,4] = ""
$data[42,5] = ""
$data[42,6] = ""
$data[43,0] = 43
$data[43,1] = 'MADRI-MADR-H-03-COS'
$data[43,2] = '21:1'
$data[43,3] = 'unlocked'
$data[43,4] = ""
$data[43,5] = ""
$data[43,6] = ""
$data[44,0] = 44
$data[44,1] = 'MADRI-MADR-H-03-COS'
$data[44,2] = '22:0'
$data[44,3] = 'unlocked'
$data[44,4] = ""
$data[44,5] = ""
$data[44,6] = ""
$data[45,0] = 45
$data[45,1] = 'MADRI-MADR-H-03-COS'
$data[45,2] = '22:1'
$data[45,3] = 'unlocked'
$data[45,4] = ""
$data[45,5] = ""
$data[45,6] = ""
$data[46,0] = 46
$data[46,1] = 'MADRI-MADR-H-03-COS'
$data[46,2] = '23:0'
$data[46,3] = 'unlocked'
$data[46,4] = ""
$data[46,5] = ""
$data[46,6] = ""
$data[47,0] = 47
$data[47,1] = 'MADRI-MADR-H-03-COS'
$data[47,2] = '23:1'
$data[47,3] = 'unlocked'
$data[47,4] = ""
$data[47,5] = ""
$data[47,6] = ""
$data[48,0] = 48
$data[48,1] = 'MADRI-MADR-H-03-COS'
$data[48,2] = '24:0'
$data[48,3] = 'unlocked'
$data[48,4] = ""
$data[48,5] = ""
$data[48,6] = ""
$data[49,0] = 49
$data[49,1] = 'MADRI-MADR-H-03-COS'
$data[49,2] = '24:1'
$data[49,3] = 'unlocked'
$data[49,4] = ""
$data[49,5] = ""
$data[49,6] = ""
$data[50,0] = 50
$data[50,1] = 'MADRI-MADR-H-03-COS'
$data[50,2] = '25:0'
$data[50,3] = 'unlocked'
$data[50,4] = ""
$data[50,5] = ""
$data[50,6] = ""
$data[51,0] = 51
$data[51,1] = 'MADRI-MADR-H-03-COS'
$data[51,2] = '25:1'
$data[51,3] = 'unlocked'
$data[51,4] = ""
$data[51,5] = ""
$data[51,6] = ""
$data[52,0] = 52
$data[52,1] = 'MADRI-MADR-H-03-COS'
$data[52,2] = '26:0'
$data[52,3] = 'unlocked'
$data[52,4] = ""
$data[52,5] = ""
$data[52,6] = ""
$data[53,0] = 53
$data[53,1] = 'MADRI-MADR-H-03-COS'
$data[53,2] = '26:1'
$data[53,3] = 'unlocked'
$data[53,4] = ""
$data[53,5] = ""
$data[53,6] = ""
$data[54,0] = 54
$data[54,1] = 'MADRI-MADR-H-03-COS'
$data[54,2] = '27:0'
$data[54,3] = 'unlocked'
$data[54,4] = ""
$data[54,5] = ""
$data[54,6] = ""
$data[55,0] = 55
$data[55,1] = 'MADRI-MADR-H-03-COS'
$data[55,2] = '27:1'
$data[55,3] = 'unlocked'
$data[55,4] = ""
$data[55,5] = ""
$data[55,6] = ""
$data[56,0] = 56
$data[56,1] = 'MADRI-MADR-H-03-COS'
$data[56,2] = '28:0'
$data[56,3] = 'unlocked'
$data[56,4] = ""
$data[56,5] = ""
$data[56,6] = ""
$data[57,0] = 57
$data[57,1] = 'MADRI-MADR-H-03-COS'
$data[57,2] = '28:1'
$data[57,3] = 'unlocked'
$data[57,4] = ""
$data[57,5] = ""
$data[57,6] = ""
$data[58,0] = 58
$data[58,1] = 'MADRI-MADR-H-03-COS'
$data[58,2] = '29:0'
$data[58,3] = 'unlocked'
$data[58,4] = ""
$data[58,5] = ""
$data[58,6] = ""
$data[59,0] = 59
$data[59,1] = 'MADRI-MADR-H-03-COS'
$data[59,2] = '29:1'
$data[59,3] = 'unlocked'
$data[59,4] = ""
$data[59,5] = ""
$data[59,6] = ""
$data[60,0] = 60
$data[60,1] = 'MADRI-MADR-H-03-COS'
$data[60,2] = '30:0'
$data[60,3] = 'unlocked'
$data[60,4] = ""
$data[60,5] = ""
$data[60,6] = ""
$data[61,0] = 61
$data[61,1] = 'MADRI-MADR-H-03-COS'
$data[61,2] = '30:1'
$data[61,3] = 'unlocked'
$data[61,4] = ""
$data[61,5] = ""
$data[61,6] = ""
$data[62,0] = 62
$data[62,1] = 'MADRI-MADR-H-03-COS'
$data[62,2] = '31:0'
$data[62,3] = 'unlocked'
$data[62,4] = ""
$data[62,5] = ""
$data[62,6] = ""
$data[63,0] = 63
$data[63,1] = 'MADRI-MADR-H-03-COS'
$data[63,2] = '31:1'
$data[63,3] = 'unlocked'
$data[63,4] = ""
$data[63,5] = ""
$data[63,6] = ""
$data[64,0] = 64
$data[64,1] = 'MADRI-MADR-H-03-COS'
$data[64,2] = '32:0'
$data[64,3] = 'unlocked'
$data[64,4] = ""
$data[64,5] = ""
$data[64,6] = ""
$data[65,0] = 65
$data[65,1] = 'MADRI-MADR-H-03-COS'
$data[65,2] = '32:1'
$data[65,3] = 'unlocked'
$data[65,4] = ""
$data[65,5] = ""
$data[65,6] = ""
$data[66,0] = 66
$data[66,1] = 'MADRI-MADR-H-03-COS'
$data[66,2] = '33:0'
$data[66,3] = 'unlocked'
$data[66,4] = ""
$data[66,5] = ""
$data[66,6] = ""
$data[67,0] = 67
$data[67,1] = 'MADRI-MADR-H-03-COS'
$data[67,2] = '33:1'
$data[67,3] = 'unlocked'
$data[67,4] = ""
$data[67,5] = ""
$data[67,6] = ""
$data[68,0] = 68
$data[68,1] = 'MADRI-MADR-H-03-COS'
$data[68,2] = '34:0'
$data[68,3] = 'unlocked'
$data[68,4] = ""
$data[68,5] = ""
$data[68,6] = ""
$data[69,0] = 69
$data[69,1] = 'MADRI-MADR-H-03-COS'
$data[69,2] = '34:1'
$data[69,3] = 'unlocked'
$data[69,4] = ""
$data[69,5] = ""
$data[69,6] = ""
$data[70,0] = 70
$data[70,1] = 'MADRI-MADR-H-03-COS'
$data[70,2] = '35:0'
$data[70,3] = 'unlocked'
$data[70,4] = ""
$data[70,5] = ""
$data[70,6] = ""
$data[71,0] = 71
$data[71,1] = 'MADRI-MADR-H-03-COS'
$data[71,2] = '35:1'
$data[71,3] = 'unlocked'
$data[71,4] = ""
$data[71,5] = ""
$data[71,6] = ""
$data[72,0] = 72
$data[72,1] = 'MADRI-MADR-H-03-COS'
$data[72,2] = '36:0'
$data[72,3] = 'unlocked'
$data[72,4] = ""
$data[72,5] = ""
$data[72,6] = ""
$data[73,0] = 73
$data[73,1] = 'MADRI-MADR-H-03-COS'
$data[73,2] = '36:1'
$data[73,3] = 'unlocked'
$data[73,4] = ""
$data[73,5] = ""
$data[73,6] = ""
$data[74,0] = 74
$data[74,1] = 'MADRI-MADR-H-03-COS'
$data[74,2] = '37:0'
$data[74,3] = 'unlocked'
$data[74,4] = ""
$data[74,5] = ""
$data[74,6] = ""
$data[75,0] = 75
$data[75,1] = 'MADRI-MADR-H-03-COS'
$data[75,2] = '37:1'
$data[75,3] = 'unlocked'
$data[75,4] = ""
$data[75,5] = ""
$data[75,6] = ""
$data[76,0] = 76
$data[76,1] = 'MADRI-MADR-H-03-COS'
$data[76,2] = '38:0'
$data[76,3] = 'unlocked'
$data[76,4] = ""
$data[76,5] = ""
$data[76,6] = ""
$data[77,0] = 77
$data[77,1] = 'MADRI-MADR-H-03-COS'
$data[77,2] = '38:1'
$data[77,3] = 'unlocked'
$data[77,4] = ""
$data[77,5] = ""
$data[77,6] = ""
$data[78,0] = 78
$data[78,1] = 'MADRI-MADR-H-03-COS'
$data[78,2] = '39:0'
$data[78,3] = 'unlocked'
$data[78,4] = ""
$data[78,5] = ""
$data[78,6] = ""
$data[79,0] = 79
$data[79,1] = 'MADRI-MADR-H-03-COS'
$data[79,2] = '39:1'
$data[79,3] = 'unlocked'
$data[79,4] = ""
$data[79,5] = ""
$data[79,6] = ""
$data[80,0] = 80
$data[80,1] = 'MADRI-MADR-H-03-COS'
$data[80,2] = '3:0'
$data[80,3] = 'unlocked'
$data[80,4] = ""
$data[80,5] = ""
$data[80,6] = ""
$data[81,0] = 81
$data[81,1] = 'MADRI-MADR-H-03-COS'
$data[81,2] = '3:1'
$data[81,3] = 'unlocked'
$data[81,4] = ""
$data[81,5] = ""
$data[81,6] = ""
$data[82,0] = 82
$data[82,1] = 'MADRI-MADR-H-03-COS'
$data[82,2] = '40:0'
$data[82,3] = 'unlocked'
$data[82,4] = ""
$data[82,5] = ""
$data[82,6] = ""
$data[83,0] = 83
$data[83,1] = 'MADRI-MADR-H-03-COS'
$data[83,2] = '40:1'
$data[83,3] = 'unlocked'
$data[83,4] = ""
$data[83,5] = ""
$data[83,6] = ""
$data[84,0] = 84
$data[84,1] = 'MADRI-MADR-H-03-COS'
$data[84,2] = '41:0'
$data[84,3] = 'unlocked'
$data[84,4] = ""
$data[84,5] = ""
$data[84,6] = ""
$data[85,0] = 85
$data[85,1] = 'MADRI-MADR-H-03-COS'
$data[85,2] = '41:1'
$data[85,3] = 'unlocked'
$data[85,4] = ""
$data[85,5] = ""
$data[85,6] = ""
$data[86,0] = 86
$data[86,1] = 'MADRI-MADR-H-03-COS'
$data[86,2] = '42:0'
$data[86,3] = 'unlocked'
$data[86,4] = ""
$data[86,5] = ""
$data[86,6] = ""
$data[87,0] = 87
$data[87,1] = 'MADRI-MADR-H-03-COS'
$data[87,2] = '42:1'
$data[87,3] = 'unlocked'
$data[87,4] = ""
$data[87,5] = ""
$data[87,6] = ""
$data[88,0] = 88
$data[88,1] = 'MADRI-MADR-H-03-COS'
$data[88,2] = '43:0'
$data[88,3] = 'unlocked'
$data[88,4] = ""
$data[88,5] = ""
$data[88,6] = ""
$data[89,0] = 89
$data[89,1] = 'MADRI-MADR-H-03-COS'
$data[89,2] = '43:1'
$data[89,3] = 'unlocked'
$data[89,4] = ""
$data[89,5] = ""
$data[89,6] = ""
$data[90,0] = 90
$data[90,1] = 'MADRI-MADR-H-03-COS'
$data[90,2] = '44:0'
$data[90,3] = 'unlocked'
$data[90,4] = ""
$data[90,5] = ""
$data[90,6] = ""
$data[91,0] = 91
$data[91,1] = 'MADRI-MADR-H-03-COS'
$data[91,2] = '44:1'
$data[91,3] = 'unlocked'
$data[91,4] = ""
$data[91,5] = ""
$data[91,6] = ""
$data[92,0] = 92
$data[92,1] = 'MADRI-MADR-H-03-COS'
$data[92,2] = '45:0'
$data[92,3] = 'unlocked'
$data[92,4] = ""
$data[92,5] = ""
$data[92,6] = ""
$data[93,0] = 93
$data[93,1] = 'MADRI-MADR-H-03-COS'
$data[93,2] = '45:1'
$data[93,3] = 'unlocked'
$data[93,4] = ""
$data[93,5] = ""
$data[93,6] = ""
$data[94,0] = 94
$data[94,1] = 'MADRI-MADR-H-03-COS'
$data[94,2] = '46:0'
$data[94,3] = 'unlocked'
$data[94,4] = ""
$data[94,5] = ""
$data[94,6] = ""
$data[95,0] = 95
$data[95,1] = 'MADRI-MADR-H-03-COS'
$data[95,2] = '46:1'
$data[95,3] = 'unlocked'
$data[95,4] = ""
$data[95,5] = ""
$data[95,6] = ""
$data[96,0] = 96
$data[96,1] = 'MADRI-MADR-H-03-COS'
$data[96,2] = '47:0'
$data[96,3] = 'unlocked'
$data[96,4] = ""
$data[96,5] = ""
$data[96,6] = ""
$data[97,0] = 97
$data[97,1] = 'MADRI-MADR-H-03-COS'
$data[97,2] = '47:1'
$data[97,3] = 'unlocked'
$data[97,4] = ""
$data[97,5] = ""
$data[97,6] = ""
$data[98,0] = 98
$data[98,1] = 'MADRI-MADR-H-03-COS'
$data[98,2] = '48:0'
$data[98,3] = 'unlocked'
$data[98,4] = ""
$data[98,5] = ""
$data[98,6] = ""
$data[99,0] = 99
$data[99,1] = 'MADRI-MADR-H-03-COS'
$data[99,2] = '48:1'
$data[99,3] = 'unlocked'
$data[99,4] = ""
$data[99,5] = ""
$data[99,6] = ""
$data[100,0] = 100
$data[100,1] = 'MADRI-MADR-H-03-COS'
$data[100,2] = '49:0'
$data[100,3] = 'unlocked'
$data[100,4] = ""
$data[100,5] = ""
$data[100,6] = ""
$data[101,0] = 101
$data[101,1] = 'MADRI-MADR-H-03-COS'
$data[101,2] = '49:1'
$data[101,3] = 'unlocked'
$data[101,4] = ""
$data[101,5] = ""
$data[101,6] = ""
$data[102,0] = 102
$data[102,1] = 'MADRI-MADR-H-03-COS'
$data[102,2] = '4:0'
$data[102,3] = 'unlocked'
$data[102,4] = ""
$data[102,5] = ""
$data[102,6] = ""
$data[103,0] = 103
$data[103,1] = 'MADRI-MADR-H-03-COS'
$data[103,2] = '4:1'
$data[103,3] = 'unlocked'
$data[103,4] = ""
$data[103,5] = ""
$data[103,6] = ""
$data[104,0] = 104
$data[104,1] = 'MADRI-MADR-H-03-COS'
$data[104,2] = '50:0'
$data[104,3] = 'unlocked'
$data[104,4] = ""
$data[104,5] = ""
$data[104,6] = ""
$data[105,0] = 105
$data[105,1] = 'MADRI-MADR-H-03-COS'
$data[105,2] = '50:1'
$data[105,3] = 'unlocked'
$data[105,4] = ""
$data[105,5] = ""
$data[105,6] = ""
$data[106,0] = 106
$data[106,1] = 'MADRI-MADR-H-03-COS'
$data[106,2] = '51:0'
$data[106,3] = 'unlocked'
$data[106,4] = ""
$data[106,5] = ""
$data[106,6] = ""
$data[107,0] = 107
$data[107,1] = 'MADRI-MADR-H-03-COS'
$data[107,2] = '51:1'
$data[107,3] = 'unlocked'
$data[107,4] = ""
$data[107,5] = ""
$data[107,6] = ""
$data[108,0] = 108
$data[108,1] = 'MADRI-MADR-H-03-COS'
$data[108,2] = '52:0'
$data[108,3] = 'unlocked'
$data[108,4] = ""
$data[108,5] = ""
$data[108,6] = ""
$data[109,0] = 109
$data[109,1] = 'MADRI-MADR-H-03-COS'
$data[109,2] = '52:1'
$data[109,3] = 'unlocked'
$data[109,4] = ""
$data[109,5] = ""
$data[109,6] = ""
$data[110,0] = 110
$data[110,1] = 'MADRI-MADR-H-03-COS'
$data[110,2] = '53:0'
$data[110,3] = 'unlocked'
$data[110,4] = ""
$data[110,5] = ""
$data[110,6] = ""
$data[111,0] = 111
$data[111,1] = 'MADRI-MADR-H-03-COS'
$data[111,2] = '53:1'
$data[111,3] = 'unlocked'
$data[111,4] = ""
$data[111,5] = ""
$data[111,6] = ""
$data[112,0] = 112
$data[112,1] = 'MADRI-MADR-H-03-COS'
$data[112,2] = '54:0'
$data[112,3] = 'unlocked'
$data[112,4] = ""
$data[112,5] = ""
$data[112,6] = ""
$data[113,0] = 113
$data[113,1] = 'MADRI-MADR-H-03-COS'
$data[113,2] = '54:1'
$data[113,3] = 'unlocked'
$data[113,4] = ""
$data[113,5] = ""
$data[113,6] = ""
$data[114,0] = 114
$data[114,1] = 'MADRI-MADR-H-03-COS'
$data[114,2] = '55:0'
$data[114,3] = 'unlocked'
$data[114,4] = ""
$data[114,5] = ""
$data[114,6] = ""
$data[115,0] = 115
$data[115,1] = 'MADRI-MADR-H-03-COS'
$data[115,2] = '55:1'
$data[115,3] = 'unlocked'
$data[115,4] = ""
$data[115,5] = ""
$data[115,6] = ""
$data[116,0] = 116
$data[116,1] = 'MADRI-MADR-H-03-COS'
$data[116,2] = '56:0'
$data[116,3] = 'unlocked'
$data[116,4] = ""
$data[116,5] = ""
$data[116,6] = ""
$data[117,0] = 117
$data[117,1] = 'MADRI-MADR-H-03-COS'
$data[117,2] = '56:1'
$data[117,3] = 'unlocked'
$data[117,4] = ""
$data[117,5] = ""
$data[117,6] = ""
$data[118,0] = 118
$data[118,1] = 'MADRI-MADR-H-03-COS'
$data[118,2] = '57:0'
$data[118,3] = 'unlocked'
$data[118,4] = ""
$data[118,5] = ""
$data[118,6] = ""
$data[119,0] = 119
$data[119,1] = 'MADRI-MADR-H-03-COS'
$data[119,2] = '57:1'
$data[119,3] = 'unlocked'
$data[119,4] = ""
$data[119,5] = ""
$data[119,6] = ""
$data[120,0] = 120
$data[120,1] = 'MADRI-MADR-H-03-COS'
$data[120,2] = '58:0'
$data[120,3] = 'unlocked'
$data[120,4] = ""
$data[120,5] = ""
$data[120,6] = ""
$data[121,0] = 121
$data[121,1] = 'MADRI-MADR-H-03-COS'
$data[121,2] = '58:1'
$data[121,3] = 'unlocked'
$data[121,4] = ""
$data[121,5] = ""
$data[121,6] = ""
$data[122,0] = 122
$data[122,1] = 'MADRI-MADR-H-03-COS'
$data[122,2] = '59:0'
$data[122,3] = 'unlocked'
$data[122,4] = ""
$data[122,5] = ""
$data[122,6] = ""
$data[123,0] = 123
$data[123,1] = 'MADRI-MADR-H-03-COS'
$data[123,2] = '59:1'
$data[123,3] = 'unlocked'
$data[123,4] = ""
$data[123,5] = ""
$data[123,6] = ""
$data[124,0] = 124
$data[124,1] = 'MADRI-MADR-H-03-COS'
$data[124,2] = '60:0'
$data[124,3] = 'unlocked'
$data[124,4] = ""
$data[124,5] = ""
$data[124,6] = ""
$data[125,0] = 125
$data[125,1] = 'MADRI-MADR-H-03-COS'
$data[125,2] = '60:1'
$data[125,3] = 'unlocked'
$data[125,4] = ""
$data[125,5] = ""
$data[125,6] = ""
$data[126,0] = 126
$data[126,1] = 'MADRI-MADR-H-03-COS'
$data[126,2] = '61:0'
$data[126,3] = 'unlocked'
$data[126,4] = ""
$data[126,5] = ""
$data[126,6] = ""
$data[127,0] = 127
$data[127,1] = 'MADRI-MADR-H-03-COS'
$data[127,2] = '61:1'
$data[127,3] = 'unlocked'
$data[127,4] = ""
$data[127,5] = ""
$data[127,6] = ""
$data[128,0] = 128
$data[128,1] = 'MADRI-MADR-H-03-COS'
$data[128,2] = '62:0'
$data[128,3] = 'unlocked'
$data[128,4] = ""
$data[128,5] = ""
$data[128,6] = ""
$data[129,0] = 129
$data[129,1] = 'MADRI-MADR-H-03-COS'
$data[129,2] = '62:1'
$data[129,3] = 'unlocked'
$data[129,4] = ""
$data[129,5] = ""
$data[129,6] = ""
$data[130,0] = 130
$data[130,1] = 'MADRI-MADR-H-03-COS'
$data[130,2] = '63:0'
$data[130,3] = 'unlocked'
$data[130,4] = ""
$data[130,5] = ""
$data[130,6] = ""
$data[131,0] = 131
$data[131,1] = 'MADRI-MADR-H-03-COS'
$data[131,2] = '63:1'
$data[131,3] = 'unlocked'
$data[131,4] = ""
$data[131,5] = ""
$data[131,6] = ""
$data[132,0] = 132
$data[132,1] = 'MADRI-MADR-H-03-COS'
$data[132,2] = '64:0'
$data[132,3] = 'unlocked'
$data[132,4] = ""
$data[132,5] = ""
$data[132,6] = ""
$data[133,0] = 133
$data[133,1] = 'MADRI-MADR-H-03-COS'
$data[133,2] = '64:1'
$data[133,3] = 'unlocked'
$data[133,4] = ""
$data[133,5] = ""
$data[133,6] = ""
$data[134,0] = 134
$data[134,1] = 'MADRI-MADR-H-03-COS'
$data[134,2] = '65:0'
$data[134,3] = 'unlocked'
$data[134,4] = ""
$data[134,5] = ""
$data[134,6] = ""
$data[135,0] = 135
$data[135,1] = 'MADRI-MADR-H-03-COS'
$data[135,2] = '65:1'
$data[135,3] = 'unlocked'
$data[135,4] = ""
$data[135,5] = ""
$data[135,6] = ""
$data[136,0] = 136
$data[136,1] = 'MADRI-MADR-H-03-COS'
$data[136,2] = '66:0'
$data[136,3] = 'unlocked'
$data[136,4] = ""
$data[136,5] = ""
$data[136,6] = ""
$data[137,0] = 137
$data[137,1] = 'MADRI-MADR-H-03-COS'
$data[137,2] = '66:1'
$data[137,3] = 'unlocked'
$data[137,4] = ""
$data[137,5] = ""
$data[137,6] = ""
$data[138,0] = 138
$data[138,1] = 'MADRI-MADR-H-03-COS'
$data[138,2] = '67:0'
$data[138,3] = 'unlocked'
$data[138,4] = ""
$data[138,5] = ""
$data[138,6] = ""
$data[139,0] = 139
$data[139,1] = 'MADRI-MADR-H-03-COS'
$data[139,2] = '67:1'
$data[139,3] = 'unlocked'
$data[139,4] = ""
$data[139,5] = ""
$data[139,6] = ""
$data[140,0] = 140
$data[140,1] = 'MADRI-MADR-H-03-COS'
$data[140,2] = '68:0'
$data[140,3] = 'unlocked'
$data[140,4] = ""
$data[140,5] = ""
$data[140,6] = ""
$data[141,0] = 141
$data[141,1] = 'MADRI-MADR-H-03-COS'
$data[141,2] = '68:1'
$data[141,3] = 'unlocked'
$data[141,4] = ""
$data[141,5] = ""
$data[141,6] = ""
$data[142,0] = 142
$data[142,1] = 'MADRI-MADR-H-03-COS'
$data[142,2] = '69:0'
$data[142,3] = 'unlocked'
$data[142,4] = ""
$data[142,5] = ""
$data[142,6] = ""
$data[143,0] = 143
$data[143,1] = 'MADRI-MADR-H-03-COS'
$data[143,2] = '69:1'
$data[143,3] = 'unlocked'
$data[143,4] = ""
$data[143,5] = ""
$data[143,6] = ""
$data[144,0] = 144
$data[144,1] = 'MADRI-MADR-H-03-COS'
$data[144,2] = '6:0'
$data[144,3] = 'unlocked'
$data[144,4] = ""
$data[144,5] = ""
$data[144,6] = ""
$data[145,0] = 145
$data[145,1] = 'MADRI-MADR-H-03-COS'
$data[145,2] = '6:1'
$data[145,3] = 'unlocked'
$data[145,4] = ""
$data[145,5] = ""
$data[145,6] = ""
$data[146,0] = 146
$data[146,1] = 'MADRI-MADR-H-03-COS'
$data[146,2] = '70:0'
$data[146,3] = 'unlocked'
$data[146,4] = ""
$data[146,5] = ""
$data[146,6] = ""
$data[147,0] = 147
$data[147,1] = 'MADRI-MADR-H-03-COS'
$data[147,2] = '70:1'
$data[147,3] = 'unlocked'
$data[147,4] = ""
$data[147,5] = ""
$data[147,6] = ""
$data[148,0] = 148
$data[148,1] = 'MADRI-MADR-H-03-COS'
$data[148,2] = '71:0'
$data[148,3] = 'unlocked'
$data[148,4] = ""
$data[148,5] = ""
$data[148,6] = ""
$data[149,0] = 149
$data[149,1] = 'MADRI-MADR-H-03-COS'
$data[149,2] = '71:1'
$data[149,3] = 'unlocked'
$data[149,4] = ""
$data[149,5] = ""
$data[149,6] = ""
$data[150,0] = 150
$data[150,1] = 'MADRI-MADR-H-03-COS'
$data[150,2] = '72:0'
$data[150,3] = 'unlocked'
$data[150,4] = ""
$data[150,5] = ""
$data[150,6] = ""
$data[151,0] = 151
$data[151,1] = 'MADRI-MADR-H-03-COS'
$data[151,2] = '72:1'
$data[151,3] = 'unlocked'
$data[151,4] = ""
$data[151,5] = ""
$data[151,6] = ""
$data[152,0] = 152
$data[152,1] = 'MADRI-MADR-H-03-COS'
$data[152,2] = '73:0'
$data[152,3] = 'unlocked'
$data[152,4] = ""
$data[152,5] = ""
$data[152,6] = ""
$data[153,0] = 153
$data[153,1] = 'MADRI-MADR-H-03-COS'
$data[153,2] = '73:1'
$data[153,3] = 'unlocked'
$data[153,4] = ""
$data[153,5] = ""
$data[153,6] = ""
$data[154,0] = 154
$data[154,1] = 'MADRI-MADR-H-03-COS'
$data[154,2] = '74:0'
$data[154,3] = 'unlocked'
$data[154,4] = ""
$data[154,5] = ""
$data[154,6] = ""
$data[155,0] = 155
$data[155,1] = 'MADRI-MADR-H-03-COS'
$data[155,2] = '74:1'
$data[155,3] = 'unlocked'
$data[155,4] = ""
$data[155,5] = ""
$data[155,6] = ""
$data[156,0] = 156
$data[156,1] = 'MADRI-MADR-H-03-COS'
$data[156,2] = '75:0'
$data[156,3] = 'unlocked'
$data[156,4] = ""
$data[156,5] = ""
$data[156,6] = ""
$data[157,0] = 157
$data[157,1] = 'MADRI-MADR-H-03-COS'
$data[157,2] = '75:1'
$data[157,3] = 'unlocked'
$data[157,4] = ""
$data[157,5] = ""
$data[157,6] = ""
$data[158,0] = 158
$data[158,1] = 'MADRI-MADR-H-03-COS'
$data[158,2] = '76:0'
$data[158,3] = 'unlocked'
$data[158,4] = ""
$data[158,5] = ""
$data[158,6] = ""
$data[159,0] = 159
$data[159,1] = 'MADRI-MADR-H-03-COS'
$data[159,2] = '76:1'
$data[159,3] = 'unlocked'
$data[159,4] = ""
$data[159,5] = ""
$data[159,6] = ""
$data[160,0] = 160
$data[160,1] = 'MADRI-MADR-H-03-COS'
$data[160,2] = '77:0'
$data[160,3] = 'unlocked'
$data[160,4] = ""
$data[160,5] = ""
$data[160,6] = ""
$data[161,0] = 161
$data[161,1] = 'MADRI-MADR-H-03-COS'
$data[161,2] = '77:1'
$data[161,3] = 'unlocked'
$data[161,4] = ""
$data[161,5] = ""
$data[161,6] = ""
$data[162,0] = 162
$data[162,1] = 'MADRI-MADR-H-03-COS'
$data[162,2] = '78:0'
$data[162,3] = 'unlocked'
$data[162,4] = ""
$data[162,5] = ""
$data[162,6] = ""
$data[163,0] = 163
$data[163,1] = 'MADRI-MADR-H-03-COS'
$data[163,2] = '78:1'
$data[163,3] = 'unlocked'
$data[163,4] = ""
$data[163,5] = ""
$data[163,6] = ""
$data[164,0] = 164
$data[164,1] = 'MADRI-MADR-H-03-COS'
$data[164,2] = '79:0'
$data[164,3] = 'unlocked'
$data[164,4] = ""
$data[164,5] = ""
$data[164,6] = ""
$data[165,0] = 165
$data[165,1] = 'MADRI-MADR-H-03-COS'
$data[165,2] = '79:1'
$data[165,3] = 'unlocked'
$data[165,4] = ""
$data[165,5] = ""
$data[165,6] = ""
$data[166,0] = 166
$data[166,1] = 'MADRI-MADR-H-03-COS'
$data[166,2] = '80:0'
$data[166,3] = 'unlocked'
$data[166,4] = ""
$data[166,5] = ""
$data[166,6] = ""
$data[167,0] = 167
$data[167,1] = 'MADRI-MADR-H-03-COS'
$data[167,2] = '80:1'
$data[167,3] = 'unlocked'
$data[167,4] = ""
$data[167,5] = ""
$data[167,6] = ""
$data[168,0] = 168
$data[168,1] = 'MADRI-MADR-H-03-COS'
$data[168,2] = '81:0'
$data[168,3] = 'unlocked'
$data[168,4] = ""
$data[168,5] = ""
$data[168,6] = ""
$data[169,0] = 169
$data[169,1] = 'MADRI-MADR-H-03-COS'
$data[169,2] = '81:1'
$data[169,3] = 'unlocked'
$data[169,4] = ""
$data[169,5] = ""
$data[169,6] = ""
$data[170,0] = 170
$data[170,1] = 'MADRI-MADR-H-03-COS'
$data[170,2] = '82:0'
$data[170,3] = 'unlocked'
$data[170,4] = ""
$data[170,5] = ""
$data[170,6] = ""
$data[171,0] = 171
$data[171,1] = 'MADRI-MADR-H-03-COS'
$data[171,2] = '82:1'
$data[171,3] = 'unlocked'
$data[171,4] = ""
$data[171,5] = ""
$data[171,6] = ""
$data[172,0] = 172
$data[172,1] = 'MADRI-MADR-H-03-COS'
$data[172,2] = '83:0'
$data[172,3] = 'unlocked'
$data[172,4] = ""
$data[172,5] = ""
$data[172,6] = ""
$data[173,0] = 173
$data[173,1] = 'MADRI-MADR-H-03-COS'
$data[173,2] = '83:1'
$data[173,3] = 'unlocked'
$data[173,4] = ""
$data[173,5] = ""
$data[173,6] = ""
$data[174,0] = 174
$data[174,1] = 'MADRI-MADR-H-03-COS'
$data[174,2] = '84:0'
$data[174,3] = 'unlocked'
$data[174,4] = ""
$data[174,5] = ""
$data[174,6] = ""
$data[175,0] = 175
$data[175,1] = 'MADRI-MADR-H-03-COS'
$data[175,2] = '84:1'
$data[175,3] = 'unlocked'
$data[175,4] = ""
$data[175,5] = ""
$data[175,6] = ""
$data[176,0] = 176
$data[176,1] = 'MADRI-MADR-H-03-COS'
$data[176,2] = '85:0'
$data[176,3] = 'unlocked'
$data[176,4] = ""
$data[176,5] = ""
$data[176,6] = ""
$data[177,0] = 177
$data[177,1] = 'MADRI-MADR-H-03-COS'
$data[177,2] = '85:1'
$data[177,3] = 'unlocked'
$data[177,4] = ""
$data[177,5] = ""
$data[177,6] = ""
$data[178,0] = 178
$data[178,1] = 'MADRI-MADR-H-03-COS'
$data[178,2] = '86:0'
$data[178,3] = 'unlocked'
$data[178,4] = ""
$data[178,5] = ""
$data[178,6] = ""
$data[179,0] = 179
$data[179,1] = 'MADRI-MADR-H-03-COS'
$data[179,2] = '86:1'
$data[179,3] = 'unlocked'
$data[179,4] = ""
$data[179,5] = ""
$data[179,6] = ""
$data[180,0] = 180
$data[180,1] = 'MADRI-MADR-H-03-COS'
$data[180,2] = '87:0'
$data[180,3] = 'unlocked'
$data[180,4] = ""
$data[180,5] = ""
$data[180,6] = ""
$data[181,0] = 181
$data[181,1] = 'MADRI-MADR-H-03-COS'
$data[181,2] = '87:1'
$data[181,3] = 'unlocked'
$data[181,4] = ""
$data[181,5] = ""
$data[181,6] = ""
$data[182,0] = 182
$data[182,1] = 'MADRI-MADR-H-03-COS'
$data[182,2] = '88:0'
$data[182,3] = 'unlocked'
$data[182,4] = ""
$data[182,5] = ""
$data[182,6] = ""
$data[183,0] = 183
$data[183,1] = 'MADRI-MADR-H-03-COS'
$data[183,2] = '88:1'
$data[183,3] = 'unlocked'
$data[183,4] = ""
$data[183,5] = ""
$data[183,6] = ""
$data[184,0] = 184
$data[184,1] = 'MADRI-MADR-H-03-COS'
$data[184,2] = '89:0'
$data[184,3] = 'unlocked'
$data[184,4] = ""
$data[184,5] = ""
$data[184,6] = ""
$data[185,0] = 185
$data[185,1] = 'MADRI-MADR-H-03-COS'
$data[185,2] = '89:1'
$data[185,3] = 'unlocked'
$data[185,4] = ""
$data[185,5] = ""
$data[185,6] = ""
$data[186,0] = 186
$data[186,1] = 'MADRI-MADR-H-03-COS'
$data[186,2] = '8:0'
$data[186,3] = 'unlocked'
$data[186,4] = ""
$data[186,5] = ""
$data[186,6] = ""
$data[187,0] = 187
$data[187,1] = 'MADRI-MADR-H-03-COS'
$data[187,2] = '8:1'
$data[187,3] = 'unlocked'
$data[187,4] = ""
$data[187,5] = ""
$data[187,6] = ""
$data[188,0] = 188
$data[188,1] = 'MADRI-MADR-H-03-COS'
$data[188,2] = '90:0'
$data[188,3] = 'unlocked'
$data[188,4] = ""
$data[188,5] = ""
$data[188,6] = ""
$data[189,0] = 189
$data[189,1] = 'MADRI-MADR-H-03-COS'
$data[189,2] = '90:1'
$data[189,3] = 'unlocked'
$data[189,4] = ""
$data[189,5] = ""
$data[189,6] = ""
$data[190,0] = 190
$data[190,1] = 'MADRI-MADR-H-03-COS'
$data[190,2] = '91:0'
$data[190,3] = 'unlocked'
$data[190,4] = ""
$data[190,5] = ""
$data[190,6] = ""
$data[191,0] = 191
$data[191,1] = 'MADRI-MADR-H-03-COS'
$data[191,2] = '91:1'
$data[191,3] = 'unlocked'
$data[191,4] = ""
$data[191,5] = ""
$data[191,6] = ""
$data[192,0] = 192
$data[192,1] = 'MADRI-MADR-H-03-COS'
$data[192,2] = '92:0'
$data[192,3] = 'unlocked'
$data[192,4] = ""
$data[192,5] = ""
$data[192,6] = ""
$data[193,0] = 193
$data[193,1] = 'MADRI-MADR-H-03-COS'
$data[193,2] = '92:1'
$data[193,3] = 'unlocked'
$data[193,4] = ""
$data[193,5] = ""
$data[193,6] = ""
$data[194,0] = 194
$data[194,1] = 'MADRI-MADR-H-03-COS'
$data[194,2] = '93:0'
$data[194,3] = 'unlocked'
$data[194,4] = ""
$data[194,5] = ""
$data[194,6] = ""
$data[195,0] = 195
$data[195,1] = 'MADRI-MADR-H-03-COS'
$data[195,2] = '93:1'
$data[195,3] = 'unlocked'
$data[195,4] = ""
$data[195,5] = ""
$data[195,6] = ""
$data[196,0] = 196
$data[196,1] = 'MADRI-MADR-H-03-COS'
$data[196,2] = '94:0'
$data[196,3] = 'unlocked'
$data[196,4] = ""
$data[196,5] = ""
$data[196,6] = ""
$data[197,0] = 197
$data[197,1] = 'MADRI-MADR-H-03-COS'
$data[197,2] = '94:1'
$data[197,3] = 'unlocked'
$data[197,4] = ""
$data[197,5] = ""
$data[197,6] = ""
$data[198,0] = 198
$data[198,1] = 'MADRI-MADR-H-03-COS'
$data[198,2] = '95:0'
$data[198,3] = 'unlocked'
$data[198,4] = ""
$data[198,5] = ""
$data[198,6] = ""
$data[199,0] = 199
$data[199,1] = 'MADRI-MADR-H-03-COS'
$data[199,2] = '95:1'
$data[199,3] = 'unlocked'
$data[199,4] = ""
$data[199,5] = ""
$data[199,6] = ""
$data[200,0] = 200
$data[200,1] = 'MADRI-MADR-H-03-COS'
$data[200,2] = '96:0'
$data[200,3] = 'unlocked'
$data[200,4] = ""
$data[200,5] = ""
$data[200,6] = ""
$data[201,0] = 201
$data[201,1] = 'MADRI-MADR-H-03-COS'
$data[201,2] = '96:1'
$data[201,3] = 'unlocked'
$data[201,4] = ""
$data[201,5] = ""
$data[201,6] = ""
$data[202,0] = 202
$data[202,1] = 'MADRI-MADR-H-03-COS'
$data[202,2] = '97:0'
$data[202,3] = 'unlocked'
$data[202,4] = ""
$data[202,5] = ""
$data[202,6] = ""
$data[203,0] = 203
$data[203,1] = 'MADRI-MADR-H-03-COS'
$data[203,2] = '97:1'
$data[203,3] = 'unlocked'
$data[203,4] = ""
$data[203,5] = ""
$data[203,6] = ""
$data[204,0] = 204
$data[204,1] = 'MADRI-MADR-H-03-COS'
$data[204,2] = '98:0'
$data[204,3] = 'unlocked'
$data[204,4] = ""
$data[204,5] = ""
$data[204,6] = ""
$data[205,0] = 205
$data[205,1] = 'MADRI-MADR-H-03-COS'
$data[205,2] = '98:1'
$data[205,3] = 'unlocked'
$data[205,4] = ""
$data[205,5] = ""
$data[205,6] = ""
$data[206,0] = 206
$data[206,1] = 'MADRI-MADR-H-03-COS'
$data[206,2] = '99:0'
$data[206,3] = 'unlocked'
$data[206,4] = ""
$data[206,5] = ""
$data[206,6] = ""
$data[207,0] = 207
$data[207,1] = 'MADRI-MADR-H-03-COS'
$data[207,2] = '99:1'
$data[207,3] = 'unlocked'
$data[207,4] = ""
$data[207,5] = ""
$data[207,6] = ""
$data[208,0] = 208
$data[208,1] = 'MADRI-MADR-H-03-COS'
$data[208,2] = '9:0'
$data[208,3] = 'unlocked'
$data[208,4] = ""
$data[208,5] = ""
$data[208,6] = ""
$data[209,0] = 209
$data[209,1] = 'MADRI-MADR-H-03-COS'
$data[209,2] = '9:1'
$data[209,3] = 'unlocked'
$data[209,4] = ""
$data[209,5] = ""
$data[209,6] = ""

$ws.Range("A2:G211").Value2 = $data

# Fix style of newly added column-A cells (rows 172-211) to match existing style (s="1")
$ws.Range("A171").Copy() | Out-Null
$ws.Range("A172:A211").PasteSpecial(-4122) | Out-Null

Write-Host "done"
